# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# Sheet "展览": row -> new F value
$exhibitionUpdates = @{
    3  = 343
    4  = 444
    5  = 1749
    6  = 90
    7  = 2214
    11 = 4998
    15 = 231
    20 = 125
    21 = 3966
    22 = 721
    23 = 696
    24 = 26
    25 = 24
    26 = 112
    27 = 123
    30 = 93
    32 = 12
    34 = 1012
    35 = 3
    36 = 2565
    38 = 23
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型": row -> new F value (row numbers shifted by +1 vs "展览"
# starting from row 35 because of an extra data row present in this sheet)
$allTypesUpdates = @{
    3  = 343
    4  = 444
    5  = 1749
    6  = 90
    7  = 2214
    11 = 4998
    15 = 231
    20 = 125
    21 = 3966
    22 = 721
    23 = 696
    24 = 26
    25 = 24
    26 = 112
    27 = 123
    30 = 93
    32 = 12
    35 = 1012
    36 = 3
    37 = 2565
    39 = 23
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
